{"js": "// Replace the three-digit-by-one-digit multiplication equations in the\n// answer table with a new set of equations (same format \"A\u00d7B=C\").\n// Each old value is unique in the document, so a simple search/replace\n// per pair is sufficient and keeps all existing run formatting intact.\nconst replacements = [\n  [\"268\u00d77=1876\", \"334\u00d72=668\"],\n  [\"571\u00d74=2284\", \"905\u00d74=3620\"],\n  [\"659\u00d72=1318\", \"342\u00d78=2736\"],\n  [\"285\u00d76=1710\", \"226\u00d76=1356\"],\n  [\"575\u00d78=4600\", \"879\u00d79=7911\"],\n  [\"841\u00d73=2523\", \"699\u00d79=6291\"],\n  [\"677\u00d76=4062\", \"606\u00d77=4242\"],\n  [\"628\u00d76=3768\", \"903\u00d73=2709\"],\n  [\"974\u00d74=3896\", \"375\u00d76=2250\"],\n  [\"190\u00d78=1520\", \"259\u00d76=1554\"],\n  [\"405\u00d74=1620\", \"857\u00d76=5142\"],\n  [\"991\u00d74=3964\", \"136\u00d77=952\"],\n  [\"476\u00d79=4284\", \"568\u00d76=3408\"],\n  [\"229\u00d74=916\", \"263\u00d75=1315\"],\n  [\"837\u00d77=5859\", \"384\u00d72=768\"],\n  [\"440\u00d76=2640\", \"473\u00d73=1419\"],\n  [\"749\u00d73=2247\", \"774\u00d74=3096\"],\n  [\"131\u00d74=524\", \"953\u00d79=8577\"],\n  [\"264\u00d77=1848\", \"775\u00d75=3875\"],\n  [\"254\u00d79=2286\", \"687\u00d76=4122\"],\n  [\"981\u00d74=3924\", \"841\u00d72=1682\"],\n  [\"438\u00d79=3942\", \"814\u00d73=2442\"],\n  [\"277\u00d73=831\", \"704\u00d74=2816\"],\n  [\"685\u00d72=1370\", \"480\u00d72=960\"],\n  [\"996\u00d77=6972\", \"612\u00d72=1224\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication equations in the\n# answer table with a new set of equations (same format \"A\u00d7B=C\").\n# Each old value is unique in the document, so Find/Replace All per pair\n# is sufficient and keeps existing run formatting intact.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find=\"268\u00d77=1876\"; Replace=\"334\u00d72=668\"},\n    @{Find=\"571\u00d74=2284\"; Replace=\"905\u00d74=3620\"},\n    @{Find=\"659\u00d72=1318\"; Replace=\"342\u00d78=2736\"},\n    @{Find=\"285\u00d76=1710\"; Replace=\"226\u00d76=1356\"},\n    @{Find=\"575\u00d78=4600\"; Replace=\"879\u00d79=7911\"},\n    @{Find=\"841\u00d73=2523\"; Replace=\"699\u00d79=6291\"},\n    @{Find=\"677\u00d76=4062\"; Replace=\"606\u00d77=4242\"},\n    @{Find=\"628\u00d76=3768\"; Replace=\"903\u00d73=2709\"},\n    @{Find=\"974\u00d74=3896\"; Replace=\"375\u00d76=2250\"},\n    @{Find=\"190\u00d78=1520\"; Replace=\"259\u00d76=1554\"},\n    @{Find=\"405\u00d74=1620\"; Replace=\"857\u00d76=5142\"},\n    @{Find=\"991\u00d74=3964\"; Replace=\"136\u00d77=952\"},\n    @{Find=\"476\u00d79=4284\"; Replace=\"568\u00d76=3408\"},\n    @{Find=\"229\u00d74=916\"; Replace=\"263\u00d75=1315\"},\n    @{Find=\"837\u00d77=5859\"; Replace=\"384\u00d72=768\"},\n    @{Find=\"440\u00d76=2640\"; Replace=\"473\u00d73=1419\"},\n    @{Find=\"749\u00d73=2247\"; Replace=\"774\u00d74=3096\"},\n    @{Find=\"131\u00d74=524\"; Replace=\"953\u00d79=8577\"},\n    @{Find=\"264\u00d77=1848\"; Replace=\"775\u00d75=3875\"},\n    @{Find=\"254\u00d79=2286\"; Replace=\"687\u00d76=4122\"},\n    @{Find=\"981\u00d74=3924\"; Replace=\"841\u00d72=1682\"},\n    @{Find=\"438\u00d79=3942\"; Replace=\"814\u00d73=2442\"},\n    @{Find=\"277\u00d73=831\"; Replace=\"704\u00d74=2816\"},\n    @{Find=\"685\u00d72=1370\"; Replace=\"480\u00d72=960\"},\n    @{Find=\"996\u00d77=6972\"; Replace=\"612\u00d72=1224\"}\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $range.Find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n}\n"}
